# edit.ps1 -- applies the "Updated cryptos list" data refresh described by the diff.
# Strategy: D column ("Price") values must remain TEXT cells (they are often formatted
# like "67.316.35" or "1.00" which Excel would otherwise auto-convert to a number).
# We force a temporary text NumberFormat, assign the literal string, then ClearFormats()
# to drop the now-unneeded style index (matching the original un-styled text cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# --- Rows 12 & 13: Cardano and Toncoin swapped position in the ranking ---
Set-TextValue "B12" "Toncoin"
Set-TextValue "C12" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D12" "5.17"
Set-TextValue "E12" "  -0.46%  "

Set-TextValue "B13" "Cardano"
Set-TextValue "C13" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D13" "0.352"
Set-TextValue "E13" "  -2.76%  "

# --- Remaining per-row Price (D) / Volume(1h) (E) updates ---
# Row 2
Set-TextValue "D2" "67.316.35"
Set-TextValue "E2" "  +0.38%  "

# Row 3
Set-TextValue "D3" "2.545.57"
Set-TextValue "E3" "  -2.60%  "

# Row 4
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  -0.03%  "

# Row 5
Set-TextValue "D5" "591.17"
Set-TextValue "E5" "  +0.23%  "

# Row 6
Set-TextValue "D6" "174.24"
Set-TextValue "E6" "  +5.46%  "

# Row 7
Set-TextValue "E7" "  +0.02%  "

# Row 8
Set-TextValue "E8" "  -0.21%  "

# Row 9
Set-TextValue "D9" "2.543.53"
Set-TextValue "E9" "  -2.68%  "

# Row 10
Set-TextValue "E10" "  +1.54%  "

# Row 11
Set-TextValue "E11" "  +1.15%  "

# Row 14
Set-TextValue "D14" "27.11"
Set-TextValue "E14" "  -0.59%  "

# Row 15
Set-TextValue "D15" "3.012.03"
Set-TextValue "E15" "  -2.50%  "

# Row 16
Set-TextValue "D16" "0.0000179"
Set-TextValue "E16" "  -0.21%  "

# Row 17
Set-TextValue "D17" "67.140.91"
Set-TextValue "E17" "  +0.18%  "

# Row 18
Set-TextValue "D18" "2.534.37"
Set-TextValue "E18" "  -2.95%  "

# Row 19
Set-TextValue "D19" "8.09"
Set-TextValue "E19" "  +3.91%  "

# Row 20
Set-TextValue "D20" "11.46"
Set-TextValue "E20" "  -2.39%  "

# Row 21
Set-TextValue "D21" "356.17"
Set-TextValue "E21" "  +0.53%  "

# Row 22
Set-TextValue "D22" "4.23"
Set-TextValue "E22" "  -1.09%  "

# Row 23
Set-TextValue "E23" "  +1.31%  "

# Row 24
Set-TextValue "D24" "2.01"
Set-TextValue "E24" "  +5.08%  "

# Row 26
Set-TextValue "D26" "70.04"
Set-TextValue "E26" "  +1.42%  "

# Row 27
Set-TextValue "D27" "10.02"
Set-TextValue "E27" "  -4.22%  "

# Row 28
Set-TextValue "D28" "2.668.43"
Set-TextValue "E28" "  -2.85%  "

# Row 29
Set-TextValue "D29" "0.999"
Set-TextValue "E29" "  -0.01%  "

# Row 30
Set-TextValue "D30" "0.0₃0999"
Set-TextValue "E30" "  +0.69%  "

# Row 31
Set-TextValue "D31" "536.72"
Set-TextValue "E31" "  -1.06%  "

# Row 32
Set-TextValue "D32" "8.28"
Set-TextValue "E32" "  +5.68%  "

# Row 33
Set-TextValue "E33" "  +1.06%  "

# Row 34
Set-TextValue "E34" "  -0.38%  "

# Row 35
Set-TextValue "E35" "  -1.26%  "

# Row 36
Set-TextValue "D36" "1.00"
Set-TextValue "E36" "  -0.02%  "

# Row 37
Set-TextValue "E37" "  -0.01%  "

# Row 38
Set-TextValue "D38" "157.50"
Set-TextValue "E38" "  -0.94%  "

# Row 39
Set-TextValue "D39" "18.83"
Set-TextValue "E39" "  -0.43%  "

# Row 40
Set-TextValue "D40" "18.45"
Set-TextValue "E40" "  +1.16%  "

# Row 41
Set-TextValue "E41" "  -1.79%  "

# Row 42
Set-TextValue "E42" "  +0.58%  "

# Row 43
Set-TextValue "D43" "5.22"
Set-TextValue "E43" "  +1.66%  "

# Row 44
Set-TextValue "E44" "  +6.99%  "

# Row 45
Set-TextValue "E45" "  +0.01%  "

# Row 46
Set-TextValue "D46" "39.87"
Set-TextValue "E46" "  -0.76%  "

# Row 47
Set-TextValue "D47" "151.38"

# Row 48
Set-TextValue "E48" "  -1.46%  "

# Row 49
Set-TextValue "D49" "0.0₆0282"
Set-TextValue "E49" "  -5.44%  "

# Row 50
Set-TextValue "D50" "3.73"
Set-TextValue "E50" "  -0.83%  "

# Row 51
Set-TextValue "D51" "1.73"
Set-TextValue "E51" "  +1.62%  "
